# Power_VRES.xlsx update: v0.0.3r -> v0.0.4r, header/label corrections
#
# - "Excl." (short code, row 4) -> "excl"
# - "long" (short code, row 4)  -> "lon"
# - Row 3 (human-readable headers): "lat" -> "Latitude", "long" -> "Longitude",
#   "YearCom" -> "Commision Year", "YearDecom" -> "Decommision Year"
# - Version label "v0.0.3r" -> "v0.0.4r"
# - F2 loses its explicit highlight style (reverts to the default "Normal" style)
#
# Applied identically to both ScenarioA and ScenarioB sheets.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("ScenarioA", "ScenarioB")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Version label in C2
    $ws.Range("C2").Value = "v0.0.4r"

    # F2 reverts to the default/unstyled cell (no explicit fill/highlight)
    $ws.Range("F2").Style = "Normal"

    # Row 4 -- short/db code for the "exclude" flag column
    $ws.Range("A4").Value = "excl"

    # Row 3 -- human readable column names
    $ws.Range("P3").Value = "Commision Year"
    $ws.Range("Q3").Value = "Decommision Year"
    $ws.Range("R3").Value = "Latitude"
    $ws.Range("S3").Value = "Longitude"

    # Row 4 -- short/db code for longitude
    $ws.Range("S4").Value = "lon"
}
